$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 3 (keeps existing style s="2" already present on A3/B3)
$ws.Range("A3").Value = "PathToPathon"
$ws.Range("B3").Value = "C:\Users\jorov\anaconda3"

# Fill in row 4
$ws.Range("A4").Value = "PathToAttachments"
$ws.Range("B4").Value = "C:\Users\jorov\OneDrive\Документы\UiPath\GazpromTest\Attachments"

# Update selection to B4
$ws.Range("B4").Select()
